# Generate Report for Handback
# Update the Correspond Handoff / Handback datetimes for the 42a69680... row
# on both the "zh-cn" and "de-de" sheets to reflect the new report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 20:55:21"
$wsZhCn.Range("H2").Value = "2016-03-12 20:55:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 20:55:24"
$wsDeDe.Range("H2").Value = "2016-03-12 20:55:41"
